$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "Gammama" label to "Gamma" (the only real content change).
$ws.Range("C3").Value = "Gamma"

# Reflect the new selection left by the author after the edit.
$ws.Range("C3").Select()
